# [GraFx Brand Kits] Initial Documentation
# Applies the content changes for the editor-comparison features sheet:
#  - tweak several ✅ / ❇️ remark markers
#  - clear / shorten a couple of footnote-style remarks in column D
#  - rename a "metadata mapping" feature row and add two new feature rows
#    under **Data Sources** (PIM metadata mapping + custom connectors)
#  - restore the selection to the area that was edited

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple in-place cell edits -------------------------------------------------

# "Support for Out-of-the-box self service" remark: drop the old note
$ws.Range("D4").Value = " "

# "Drop shadows" / "Bullets" / "Support for mixed ink colors" / "Support for gradients"
# now carry the "being developed" marker instead of blank
$ws.Range("B24").Value = "❇️"
$ws.Range("B41").Value = "❇️"
$ws.Range("B82").Value = "❇️"

# "Edit text inline" remark: drop the redundant "GraFx Studio:" prefix
$ws.Range("D33").Value = "in Template Designer Workspace"

# "Upload images from local client" / "Support for mixed ink colors" /
# "AI powered Smart Crop" / "Type: Long text" now supported
$ws.Range("B58").Value = "✅"
$ws.Range("B58").Font.Color = 0
$ws.Range("B65").Value = "✅"
$ws.Range("B81").Value = "✅"
$ws.Range("B114").Value = "✅"

# "Support for metadata mapping from DAM" -> renamed to PIM
$ws.Range("A158").Value = "Support for metadata mapping from PIM"

# --- Insert the two new feature rows under it -----------------------------------

$ws.Rows("159:160").Insert()

$ws.Range("A159").Value = "Support for metadata mapping from DAM"
$ws.Range("B159").Value = "✅"
$ws.Range("C159").Value = " "
$ws.Range("D159").Value = " "

$ws.Range("A160").Value = "Support for custom connectors (media & data)"
$ws.Range("B160").Value = "✅"
$ws.Range("C160").Value = " "
$ws.Range("D160").Value = " "

# --- Restore selection/scroll position ------------------------------------------

$ws.Range("D64").Select()
